$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for the cells involved in the rotation
$A5 = $ws.Range("A5").Value()
$Q5 = $ws.Range("Q5").Value()
$R5 = $ws.Range("R5").Value()
$AC5 = $ws.Range("AC5").Value()

$A6 = $ws.Range("A6").Value()
$Q6 = $ws.Range("Q6").Value()
$R6 = $ws.Range("R6").Value()

$A8 = $ws.Range("A8").Value()
$Q8 = $ws.Range("Q8").Value()
$R8 = $ws.Range("R8").Value()
$AC8 = $ws.Range("AC8").Value()

# Row 5 gets the values that used to be in row 8
$ws.Range("A5").Value = $A8
$ws.Range("Q5").Value = $Q8
$ws.Range("R5").Value = $R8
$ws.Range("AC5").Value = $AC8

# Row 6 gets the values that used to be in row 5
$ws.Range("A6").Value = $A5
$ws.Range("Q6").Value = $Q5
$ws.Range("R6").Value = $R5

# Row 8 gets the values that used to be in row 6, and its comment is cleared
$ws.Range("A8").Value = $A6
$ws.Range("Q8").Value = $Q6
$ws.Range("R8").Value = $R6
$ws.Range("AC8").Value = ""
